# "Generate Report for Handback"
# The handback transform failed for this item, so the report now records
# the failure status plus the underlying error detail, for both the
# zh-cn and de-de target languages (and the roll-up Overview sheet).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"
$errorDetail = "The translationStateItem 980dd4e7afbfe618b50c9f45743d5bb988b4b0b3 is not found."

# Status column: was "Ready for handoff", now "Handback transform failed".
# (Overview rolls up both language columns; zh-cn/de-de each have their own Status cell.)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Error Detail column: previously blank, now carries the failure reason.
$zhcn.Range("P2").Value = $errorDetail
$dede.Range("P2").Value = $errorDetail

# The columns holding the longer text grew to fit the new content.
$overview.Columns.Item(5).ColumnWidth = 23.75
$overview.Columns.Item(6).ColumnWidth = 23.75

$zhcn.Columns.Item(3).ColumnWidth = 23.75
$zhcn.Columns.Item(16).ColumnWidth = 39.09

$dede.Columns.Item(3).ColumnWidth = 23.75
$dede.Columns.Item(16).ColumnWidth = 39.09
